$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column range to Text format so numeric-looking
# strings (e.g. "61.673.65", "1.00", "0.0000148") are preserved exactly
# as text instead of being parsed into floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Coin name (B) updates ---
$ws.Range("B20").Value = "Chainlink"
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("B26").Value = "Kaspa"
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("B28").Value = "PEPE"
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("B41").Value = "Maker"
$ws.Range("B42").Value = "OKB"

# --- Link (C) updates ---
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

# --- Price (D) updates ---
$ws.Range("D2").Value = "61.673.65"
$ws.Range("D3").Value = "2.979.96"
$ws.Range("D5").Value = "544.12"
$ws.Range("D6").Value = "153.07"
$ws.Range("D8").Value = "0.576"
$ws.Range("D9").Value = "2.987.91"
$ws.Range("D10").Value = "0.115"
$ws.Range("D13").Value = "3.499.81"
$ws.Range("D15").Value = "61.713.93"
$ws.Range("D16").Value = "23.82"
$ws.Range("D17").Value = "2.978.45"
$ws.Range("D18").Value = "0.0000148"
$ws.Range("D20").Value = "12.08"
$ws.Range("D21").Value = "382.75"
$ws.Range("D22").Value = "6.70"
$ws.Range("D23").Value = "1.00"
$ws.Range("D24").Value = "65.81"
$ws.Range("D25").Value = "0.473"
$ws.Range("D26").Value = "0.191"
$ws.Range("D27").Value = "3.101.29"
$ws.Range("D28").Value = "0.0₃0951"
$ws.Range("D29").Value = "0.997"
$ws.Range("D30").Value = "8.32"
$ws.Range("D33").Value = "20.53"
$ws.Range("D34").Value = "160.89"
$ws.Range("D35").Value = "4.72"
$ws.Range("D36").Value = "5.97"
$ws.Range("D38").Value = "1.28"
$ws.Range("D39").Value = "1.57"
$ws.Range("D40").Value = "3.95"
$ws.Range("D41").Value = "2.424.99"
$ws.Range("D42").Value = "37.46"
$ws.Range("D43").Value = "22.28"
$ws.Range("D44").Value = "0.671"
$ws.Range("D46").Value = "5.12"
$ws.Range("D49").Value = "271.06"
$ws.Range("D50").Value = "19.94"
$ws.Range("D51").Value = "0.0958"

# Restore the original (default/unformatted) style on the Price column
# now that the text values are safely stored, so no residual cell-level
# style reference is left behind.
$priceRange.Style = "Normal"

# --- Volume(1h) (E) updates ---
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("E3").Value = "  -4.96%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("E6").Value = "  -5.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("E9").Value = "  -5.05%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -6.42%  "
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("E17").Value = "  -5.40%  "
$ws.Range("E18").Value = "  -3.70%  "
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("E21").Value = "  -4.36%  "
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  -8.75%  "
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -6.05%  "
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("E51").Value = "  -1.61%  "

